$wb = $excel.ActiveWorkbook

# --- CVS sheet: a couple of quantities changed / one price became "NA" ---
$cvs = $wb.Worksheets.Item("CVS")
$cvs.Range("C3").Value = 0
$cvs.Range("C4").Value = "NA"

# --- PetSmart sheet: a price became "NA" / another price reset to 0 ---
$pet = $wb.Worksheets.Item("PetSmart")
$pet.Range("C2").Value = "NA"
$pet.Range("C4").Value = 0

# Leave the selection on CVS!C4 (matches the last place the user clicked there)
$cvs.Select()
$cvs.Range("C4").Select()

# Finish with PetSmart active/selected at C2 - this becomes the active tab
$pet.Select()
$pet.Range("C2").Select()
